$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 4 (Caso 1939 - SANCHEZ, MIGUEL B. AV. 1050).
# This shifts rows 5..50 up to become rows 4..49, and the sheet's used
# range shrinks from A1:N50 to A1:N49 automatically.
$ws.Rows(4).Delete()
